$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.874.78'
$ws.Range("E2").Value = '  +4.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.344.21'
$ws.Range("E3").Value = '  +2.84%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.90'
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.21'
$ws.Range("E6").Value = '  +4.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.341.68'
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("E10").Value = '  +6.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.20'
$ws.Range("E12").Value = '  +6.23%  '
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.758.51'
$ws.Range("E14").Value = '  +2.69%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.64'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.779.76'
$ws.Range("E16").Value = '  +3.87%  '
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.355.16'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.40'
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.49'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.52'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.70'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("E26").Value = '  +7.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.75'
$ws.Range("E27").Value = '  +4.00%  '
$ws.Range("E28").Value = '  +10.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.35'
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0741'
$ws.Range("E30").Value = '  +5.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.68'
$ws.Range("E31").Value = '  +3.32%  '
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.24'
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("E36").Value = '  +3.56%  '
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("E38").Value = '  +5.01%  '
$ws.Range("E39").Value = '  +7.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.75'
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("E42").Value = '  +5.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.74'
$ws.Range("E43").Value = '  +6.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '277.24'
$ws.Range("E44").Value = '  +10.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.08'
$ws.Range("E45").Value = '  +5.51%  '
$ws.Range("E46").Value = '  +2.96%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0218'
$ws.Range("E49").Value = '  +5.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.378'
$ws.Range("E51").Value = '  +2.29%  '
